$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$tr = $s.Shapes.Item(1).TextFrame.TextRange

# The title paragraph currently holds three separate runs
# ("Below", " ", "section-level"). Re-assigning the text of the whole
# range (addressed via Characters, which forces a full rewrite of the
# underlying runs instead of the minimal-diff behavior of a plain
# TextRange.Text assignment) collapses them into a single run while
# keeping the resulting visible text identical.
$full = $tr.Characters(1, $tr.Length)
$full.Text = "Below section-level"
